# Update cryptos list (prices + 1h volume %) per Oct 21 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.687.07'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '1.604.87'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''212.74'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '''0.517'
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").Value = '''28.01'
$ws.Range("E8").Value = '  +6.55%  '
$ws.Range("D9").Value = '''0.253'
$ws.Range("E9").Value = '  +1.80%  '
$ws.Range("D10").Value = '''0.0604'
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("D11").Value = '''0.0910'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '1.832.49'
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").Value = '1.604.95'
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").Value = '''0.551'
$ws.Range("E14").Value = '  +5.16%  '
$ws.Range("D15").Value = '29.679.18'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '''3.76'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '''64.04'
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").Value = '''242.63'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").Value = '''7.83'
$ws.Range("E19").Value = '  +4.78%  '
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").Value = '''4.03'
$ws.Range("E22").Value = '  +1.17%  '
$ws.Range("D23").Value = '''9.42'
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '''155.46'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '''6.45'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '''0.0483'
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").Value = '''3.20'
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("D34").Value = '1.430.79'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = '''2.94'
$ws.Range("E35").Value = '  +4.60%  '
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").Value = '''57.95'
$ws.Range("E40").Value = '  +8.41%  '
$ws.Range("D41").Value = '''0.548'
$ws.Range("E41").Value = '  +2.98%  '
$ws.Range("D42").Value = '''0.0499'
$ws.Range("E42").Value = '  +6.28%  '
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").Value = '''0.999'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = '''66.45'
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("D47").Value = '''0.980'
$ws.Range("E47").Value = '  +17.07%  '
$ws.Range("D48").Value = '''5.35'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '1.743.14'
$ws.Range("E49").Value = '  +1.26%  '
$ws.Range("D50").Value = '''86.79'
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0105'
$ws.Range("E51").Value = '  +3.36%  '

# Cells above hold numeric-looking text (e.g. "1.00", "0.0910") that must stay
# text, not auto-converted numbers; the leading apostrophe forced that, and we
# now strip the resulting quote-prefix formatting so styling matches the source.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
